$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '55.950.45'
$ws.Range("E2").Value = '  +8.59%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.217.61'
$ws.Range("E3").Value = '  +3.72%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '399.42'
$ws.Range("E5").Value = '  +3.60%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '109.69'
$ws.Range("E6").Value = '  +6.13%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.555'
$ws.Range("E7").Value = '  +2.92%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.623'
$ws.Range("E9").Value = '  +6.73%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.26'
$ws.Range("E10").Value = '  +5.84%  '
$ws.Range("E11").Value = '  +5.68%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.140'
$ws.Range("E12").Value = '  +1.91%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.733.12'
$ws.Range("E13").Value = '  +4.04%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '19.06'
$ws.Range("E14").Value = '  +2.49%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.07'
$ws.Range("E15").Value = '  +2.85%  '
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.218.71'
$ws.Range("E16").Value = '  +3.97%  '
$ws.Range("B17").Value = 'Polygon'
$ws.Range("C17").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.05'
$ws.Range("E17").Value = '  +6.02%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.59'
$ws.Range("E18").Value = '  -3.73%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '55.828.58'
$ws.Range("E19").Value = '  +8.29%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.37'
$ws.Range("E20").Value = '  +2.71%  '
$ws.Range("E21").Value = '  +6.00%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '13.08'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '302.55'
$ws.Range("E23").Value = '  +13.68%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '74.97'
$ws.Range("E24").Value = '  +7.29%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.23'
$ws.Range("E25").Value = '  +2.63%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.24'
$ws.Range("E26").Value = '  +1.83%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '28.32'
$ws.Range("E27").Value = '  +4.69%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.56'
$ws.Range("E28").Value = '  +4.07%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.172'
$ws.Range("E29").Value = '  +3.43%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("E30").Value = '  -0.07%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '11.38'
$ws.Range("E31").Value = '  +9.92%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.110'
$ws.Range("E32").Value = '  +3.20%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0495'
$ws.Range("E33").Value = '  +4.28%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '36.25'
$ws.Range("E34").Value = '  +2.96%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.11'
$ws.Range("E35").Value = '  +2.24%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '51.41'
$ws.Range("E36").Value = '  +2.22%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.52'
$ws.Range("E37").Value = '  +4.56%  '
$ws.Range("E38").Value = '  +23.53%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.999'
$ws.Range("E39").Value = '  +0.00%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '134.34'
$ws.Range("E40").Value = '  +4.22%  '
$ws.Range("B41").Value = 'NEARProtocol'
$ws.Range("C41").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.03'
$ws.Range("E41").Value = '  +9.56%  '
$ws.Range("B42").Value = 'ARBITRUM'
$ws.Range("C42").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.92'
$ws.Range("E42").Value = '  +1.62%  '
$ws.Range("B43").Value = 'Stellar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.119'
$ws.Range("E43").Value = '  +3.03%  '
$ws.Range("B44").Value = 'Celestia'
$ws.Range("C44").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '17.03'
$ws.Range("E44").Value = '  +2.86%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.283'
$ws.Range("E45").Value = '  -3.35%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '22.30'
$ws.Range("E46").Value = '  -1.13%  '
$ws.Range("E47").Value = '  +2.33%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.153.04'
$ws.Range("E48").Value = '  +4.24%  '
$ws.Range("B49").Value = 'ApeXProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.46'
$ws.Range("E49").Value = '  -1.16%  '
$ws.Range("B50").Value = 'ThetaToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.08'
$ws.Range("E50").Value = '  +41.52%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0363'
$ws.Range("E51").Value = '  +9.57%  '
